$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '26.907.61'
$ws.Range("E2").Value = '  -0.22%  '
$ws.Range("D3").Value = '1.862.50'
$ws.Range("E3").Value = '  +0.16%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9999'
$ws.Range("E4").Value = '  -0.15%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '305.03'
$ws.Range("E5").Value = '  -0.21%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9997'
$ws.Range("E6").Value = '  -0.11%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5041'
$ws.Range("E7").Value = '  -0.88%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3651'
$ws.Range("E8").Value = '  -2.28%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07181'
$ws.Range("E9").Value = '  +0.94%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8958'
$ws.Range("E10").Value = '  +1.03%  '
$ws.Range("E11").Value = '  +0.96%  '
$ws.Range("D12").Value = '1.878.35'
$ws.Range("E12").Value = '  +1.78%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.07495'
$ws.Range("E13").Value = '  -0.71%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '94.92'
$ws.Range("E14").Value = '  +6.53%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '5.235'
$ws.Range("E15").Value = '  -0.99%  '
$ws.Range("E16").Value = '  -0.16%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008505'
$ws.Range("E17").Value = '  +1.85%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '14.22'
$ws.Range("E18").Value = '  +1.10%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.9998'
$ws.Range("E19").Value = '  -0.10%  '
$ws.Range("D20").Value = '26.958.76'
$ws.Range("E20").Value = '  -0.22%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.027'
$ws.Range("E21").Value = '  -0.43%  '
$ws.Range("D22").Value = '2.114.12'
$ws.Range("E22").Value = '  +1.13%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.39'
$ws.Range("E23").Value = '  -0.88%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.412'
$ws.Range("E24").Value = '  -0.69%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '147.41'
$ws.Range("E25").Value = '  -0.86%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '1.787'
$ws.Range("E26").Value = '  -3.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.88'
$ws.Range("E27").Value = '  -0.42%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.091'
$ws.Range("E28").Value = '  +0.63%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '113.32'
$ws.Range("E29").Value = '  +0.63%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.701'
$ws.Range("E30").Value = '  +0.83%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.678'
$ws.Range("E31").Value = '  +0.77%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09227'
$ws.Range("E32").Value = '  +2.04%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05143'
$ws.Range("E33").Value = '  +0.63%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7505'
$ws.Range("E34").Value = '  +3.18%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.965'
$ws.Range("E35").Value = '  -2.97%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.154'
$ws.Range("E36").Value = '  +0.21%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.249'
$ws.Range("E37").Value = '  +6.92%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '2.588'
$ws.Range("E38").Value = '  +5.13%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.02004'
$ws.Range("E39").Value = '  -1.77%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5568'
$ws.Range("E40").Value = '  +4.71%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.069'
$ws.Range("E41").Value = '  -0.17%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '6.563'
$ws.Range("E42").Value = '  -0.02%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '116.62'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.582'
$ws.Range("E44").Value = '  +3.42%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.1473'
$ws.Range("E45").Value = '  +0.25%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.4699'
$ws.Range("E46").Value = '  +2.08%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.9993'
$ws.Range("E47").Value = '  -0.11%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.03'
$ws.Range("E48").Value = '  -0.05%  '
$ws.Range("E49").Value = '  +0.05%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '36.75'
$ws.Range("E50").Value = '  +0.73%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '63.14'
$ws.Range("E51").Value = '  -1.24%  '
